$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Convention change to support multi-axle vehicles: the sheet (and the
# corresponding shared-string labels used by the H2:H4 "class / Instance /
# Type" column) is renamed from "DrivelineAxle1_None" to "Axle1_None".
$ws.Name = "Axle1_None"

# Column B was narrowed slightly as part of the same layout tweak.
# (ColumnWidth is expressed in characters; the underlying stored width is
# quantized to 1/6-character steps by this engine, so 11.8333... is the
# value that lands on the closest achievable width to the target 12.664...)
$ws.Columns.Item(2).ColumnWidth = 11.8333333333333

# Restore the active selection on the (now frozen) bottom-right pane to H5.
$ws.Range("H5").Select() | Out-Null
